$d = $word.ActiveDocument

$replacements = @(
    @{old="86×62=5332"; new="21×88=1848"},
    @{old="60×55=3300"; new="49×93=4557"},
    @{old="40×84=3360"; new="16×99=1584"},
    @{old="19×70=1330"; new="11×58=638"},
    @{old="30×58=1740"; new="29×64=1856"},
    @{old="47×46=2162"; new="21×26=546"},
    @{old="20×72=1440"; new="53×21=1113"},
    @{old="88×46=4048"; new="69×31=2139"},
    @{old="62×17=1054"; new="16×27=432"},
    @{old="21×82=1722"; new="49×91=4459"},
    @{old="14×19=266"; new="13×98=1274"},
    @{old="31×71=2201"; new="68×50=3400"},
    @{old="25×80=2000"; new="64×61=3904"},
    @{old="58×29=1682"; new="49×94=4606"},
    @{old="57×99=5643"; new="11×72=792"},
    @{old="84×21=1764"; new="25×26=650"},
    @{old="31×60=1860"; new="37×54=1998"},
    @{old="31×79=2449"; new="72×70=5040"},
    @{old="20×54=1080"; new="37×27=999"},
    @{old="38×99=3762"; new="20×68=1360"},
    @{old="33×32=1056"; new="66×12=792"},
    @{old="26×56=1456"; new="51×96=4896"},
    @{old="67×77=5159"; new="63×78=4914"},
    @{old="88×14=1232"; new="96×49=4704"},
    @{old="25×83=2075"; new="77×85=6545"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
